$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row of data (new shared string "Bootcamp on Selenium with Java")
$ws.Range("A3").Value = "Bootcamp on Selenium with Java"

# Extend the selection to cover the new row, like the author's saved view
[void]$ws.Range("A1:A3").Select()

# Column A widened to fit the new (longer) text
$ws.Columns.Item(1).ColumnWidth = 32.6667
